$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.188.30"
$ws.Range("D3").Value = "'2.077.50"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'254.82"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "'0.682"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("D7").Value = "'62.24"
$ws.Range("E7").Value = "  +19.96%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.394"
$ws.Range("E9").Value = "  +5.06%  "
$ws.Range("D10").Value = "'61.66"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("E11").Value = "  +7.88%  "
$ws.Range("E12").Value = "  +3.25%  "
$ws.Range("E13").Value = "  +7.49%  "
$ws.Range("D14").Value = "'2.376.49"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "'0.829"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "'5.58"
$ws.Range("E16").Value = "  +9.30%  "
$ws.Range("D17").Value = "'2.077.02"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "'37.160.33"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'15.69"
$ws.Range("E19").Value = "  +12.68%  "
$ws.Range("D20").Value = "'74.93"
$ws.Range("E20").Value = "  +4.04%  "
$ws.Range("D21").Value = "'0.0₃0934"
$ws.Range("E21").Value = "  +12.10%  "
$ws.Range("D22").Value = "'5.50"
$ws.Range("E22").Value = "  +5.74%  "
$ws.Range("D23").Value = "'241.89"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'2.44"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").Value = "'2.36"
$ws.Range("E26").Value = "  +18.81%  "
$ws.Range("D27").Value = "'170.27"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'9.40"
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("D29").Value = "'20.51"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  +3.40%  "
$ws.Range("D31").Value = "'4.86"
$ws.Range("E31").Value = "  +8.39%  "
$ws.Range("D32").Value = "'1.12"
$ws.Range("E32").Value = "  +5.89%  "
$ws.Range("D33").Value = "'0.0643"
$ws.Range("E33").Value = "  +5.58%  "
$ws.Range("D34").Value = "'4.49"
$ws.Range("E34").Value = "  +9.79%  "
$ws.Range("D35").Value = "'0.0905"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D37").Value = "'2.31"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("E38").Value = "  +30.19%  "
$ws.Range("E39").Value = "  -3.18%  "
$ws.Range("D40").Value = "'1.38"
$ws.Range("E40").Value = "  +4.75%  "
$ws.Range("D41").Value = "'18.26"
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("D42").Value = "'0.0229"
$ws.Range("E42").Value = "  +2.34%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "'4.60"
$ws.Range("E43").Value = "  +32.36%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'1.17"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'99.49"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").Value = "'2.82"
$ws.Range("E46").Value = "  +3.26%  "
$ws.Range("E47").Value = "  +14.76%  "
$ws.Range("E48").Value = "  +9.33%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'1.309.73"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "'2.96"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").Value = "'6.98"
$ws.Range("E51").Value = "  -0.14%  "
